$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-15 were labeled "Right" / "Left" (river bank side); the source
# classification was simplified to a single "Sides" category.
foreach ($r in 3..15) {
    $ws.Range("D$r").Value = "Sides"
}

# Reset the view: scroll back to the top-left and move the selection to D16.
$ws.Range("D16").Select() | Out-Null
